$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Rescale existing D column values (thousands -> actual units) ---
$ws.Range("D4").Value = 9.259
$ws.Range("D5").Value = 9.375
$ws.Range("D6").Value = 9.863
$ws.Range("D7").Value = 9.796
$ws.Range("D8").Value = 9.596
$ws.Range("D9").Value = 8.896

# --- Row 14: was "Горячий ключ" 2023; becomes "Горячий ключ" 2022 (only U14 populated) ---
$ws.Range("B14").Value = 2022
$ws.Range("C14").ClearContents()
$ws.Range("F14").ClearContents()
$ws.Range("O14").ClearContents()
$ws.Range("S14").ClearContents()
$ws.Range("T14").ClearContents()
$ws.Cells.Item(14, 21).HorizontalAlignment = -4108
$ws.Range("U14").Value = 1148

# --- Row 15: "Горячий ключ" 2023 (full row, all 21 columns present) ---
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(15, $c).HorizontalAlignment = -4108
}
$ws.Range("A15").Value = "Горячий ключ"
$ws.Range("B15").Value = 2023
$ws.Range("C15").Value = 41.482
$ws.Range("D15").Value = 7.651
$ws.Range("E15").Value = 162
$ws.Range("F15").Value = 42049
$ws.Range("N15").Value = 885
$ws.Range("O15").Value = 4270.9
$ws.Range("Q15").Value = 95.557
$ws.Range("R15").Value = 1198
$ws.Range("S15").Value = 14680.8
$ws.Range("T15").Value = 210.2
$ws.Range("U15").Value = 1165

# --- Row 16: "Анапа" 2022 (full row, only U has a value) ---
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(16, $c).HorizontalAlignment = -4108
}
$ws.Range("A16").Value = "Анапа"
$ws.Range("B16").Value = 2022
$ws.Range("U16").Value = 704

# --- Row 17: "Геленджик" 2022 (full row, only U has a value) ---
for ($c = 1; $c -le 21; $c++) {
    $ws.Cells.Item(17, $c).HorizontalAlignment = -4108
}
$ws.Range("A17").Value = "Геленджик"
$ws.Range("B17").Value = 2022
$ws.Range("U17").Value = -366

# --- Row 18: "Анапа" 2023 (sparse row: A,B,D,E,F,N,O,Q,R,S,T,U) ---
foreach ($c in 1,2,4,5,6,14,15,17,18,19,20,21) {
    $ws.Cells.Item(18, $c).HorizontalAlignment = -4108
}
$ws.Range("A18").Value = "Анапа"
$ws.Range("B18").Value = 2023
$ws.Range("D18").Value = 30.729
$ws.Range("E18").Value = 643
$ws.Range("F18").Value = 49270
$ws.Range("N18").Value = 3373
$ws.Range("O18").Value = 5465.4
$ws.Range("Q18").Value = 592.065
$ws.Range("R18").Value = 7363
$ws.Range("S18").Value = 49034.2
$ws.Range("T18").Value = 3455.418
$ws.Range("U18").Value = 4899

# --- Row 19: "Геленджик" 2023 (sparse row: A,B,D,E,F,N,O,P,Q,R,S,T,U) ---
foreach ($c in 1,2,4,5,6,14,15,16,17,18,19,20,21) {
    $ws.Cells.Item(19, $c).HorizontalAlignment = -4108
}
$ws.Range("A19").Value = "Геленджик"
$ws.Range("B19").Value = 2023
$ws.Range("D19").Value = 19.093
$ws.Range("E19").Value = 96
$ws.Range("F19").Value = 56182
$ws.Range("N19").Value = 2396
$ws.Range("O19").Value = 3456.8
$ws.Range("P19").Value = 12221.3
$ws.Range("Q19").Value = 114.619
$ws.Range("R19").Value = 646
$ws.Range("S19").Value = 32472.3
$ws.Range("T19").Value = 1040.331
$ws.Range("U19").Value = 407

# --- Update the sheet view: move active selection ---
$ws.Range("J26").Select()
